$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.069.33"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.875.37"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'313.27"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "'0.5094"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").Value = "'0.3850"
$ws.Range("E8").Value = "  -2.09%  "
$ws.Range("D9").Value = "'0.09169"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").Value = "'1.123"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'41.55"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'6.341"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'20.74"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.876.37"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.204"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001115"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "'91.09"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.06593"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'18.16"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'6.110"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "28.097.20"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'11.42"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.277"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("D26").Value = "2.092.36"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.541"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.80"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'157.82"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'126.71"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "'1.067"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1053"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.612"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.599"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "'9.681"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06574"
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02429"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2179"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.209"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.263"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").Value = "'11.60"
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6400"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D43").Value = "'4.916"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.21"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6018"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "'1.235"
$ws.Range("E48").Value = "  +4.45%  "
$ws.Range("D49").Value = "'1.999"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("D50").Value = "'121.47"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("E51").Value = "  +1.50%  "
